$wb = $excel.ActiveWorkbook

# --- Sheet "10per change": E21 inlineStr "533758" -> numeric 533758, append row 22 ---
$ws1 = $wb.Worksheets.Item("10per change")

$ws1.Cells.Item(21, 5).Value = 533758

$ws1.Cells.Item(22, 1).Value = "13/06/2024 09:45:30"
$ws1.Cells.Item(22, 2).Value = 1
$ws1.Cells.Item(22, 3).Value = "APLAPOLLO"
$ws1.Cells.Item(22, 4).Value = "Apl Apollo Tubes Limited"

$c1 = $ws1.Cells.Item(22, 5)
$c1.NumberFormat = "@"
$c1.Value = "533758"
$c1.Style = "Normal"

$ws1.Cells.Item(22, 6).Value = -2.46
$ws1.Cells.Item(22, 7).Value = 1542
$ws1.Cells.Item(22, 8).Value = 593593

# --- Sheet "3 V 0.3": E7 inlineStr "532900" -> numeric 532900, append row 8 ---
$ws2 = $wb.Worksheets.Item("3 V 0.3")

$ws2.Cells.Item(7, 5).Value = 532900

$ws2.Cells.Item(8, 1).Value = "13/06/2024 09:45:30"
$ws2.Cells.Item(8, 2).Value = 1
$ws2.Cells.Item(8, 3).Value = "PAISALO"
$ws2.Cells.Item(8, 4).Value = "Paisalo Digital Ltd"

$c2 = $ws2.Cells.Item(8, 5)
$c2.NumberFormat = "@"
$c2.Value = "532900"
$c2.Style = "Normal"

$ws2.Cells.Item(8, 6).Value = 6.44
$ws2.Cells.Item(8, 7).Value = 73.42
$ws2.Cells.Item(8, 8).Value = 5265131
